$d = $word.ActiveDocument
$d.Content.Find.Execute("10+36=", $true, $false, $false, $false, $false, $true, 1, $false, "86-23=", 2) | Out-Null
$d.Content.Find.Execute("55-8=", $true, $false, $false, $false, $false, $true, 1, $false, "92-36=", 2) | Out-Null
$d.Content.Find.Execute("8+29=", $true, $false, $false, $false, $false, $true, 1, $false, "12+77=", 2) | Out-Null
$d.Content.Find.Execute("18+41=", $true, $false, $false, $false, $false, $true, 1, $false, "65+15=", 2) | Out-Null
$d.Content.Find.Execute("1+10=", $true, $false, $false, $false, $false, $true, 1, $false, "79-27=", 2) | Out-Null
$d.Content.Find.Execute("19+38=", $true, $false, $false, $false, $false, $true, 1, $false, "20+39=", 2) | Out-Null
$d.Content.Find.Execute("16+50=", $true, $false, $false, $false, $false, $true, 1, $false, "16+73=", 2) | Out-Null
$d.Content.Find.Execute("65+25=", $true, $false, $false, $false, $false, $true, 1, $false, "83+16=", 2) | Out-Null
$d.Content.Find.Execute("41+50=", $true, $false, $false, $false, $false, $true, 1, $false, "63+34=", 2) | Out-Null
$d.Content.Find.Execute("57-54=", $true, $false, $false, $false, $false, $true, 1, $false, "52+14=", 2) | Out-Null
$d.Content.Find.Execute("5+17=", $true, $false, $false, $false, $false, $true, 1, $false, "55+38=", 2) | Out-Null
$d.Content.Find.Execute("68-7=", $true, $false, $false, $false, $false, $true, 1, $false, "36+52=", 2) | Out-Null
$d.Content.Find.Execute("7+83=", $true, $false, $false, $false, $false, $true, 1, $false, "7+68=", 2) | Out-Null
$d.Content.Find.Execute("40+24=", $true, $false, $false, $false, $false, $true, 1, $false, "14+17=", 2) | Out-Null
$d.Content.Find.Execute("71-48=", $true, $false, $false, $false, $false, $true, 1, $false, "34+37=", 2) | Out-Null
$d.Content.Find.Execute("19+63=", $true, $false, $false, $false, $false, $true, 1, $false, "60+39=", 2) | Out-Null
$d.Content.Find.Execute("80-79=", $true, $false, $false, $false, $false, $true, 1, $false, "15+29=", 2) | Out-Null
$d.Content.Find.Execute("61-52=", $true, $false, $false, $false, $false, $true, 1, $false, "23+41=", 2) | Out-Null
$d.Content.Find.Execute("24+74=", $true, $false, $false, $false, $false, $true, 1, $false, "37+34=", 2) | Out-Null
$d.Content.Find.Execute("88+10=", $true, $false, $false, $false, $false, $true, 1, $false, "79-47=", 2) | Out-Null
$d.Content.Find.Execute("87-60=", $true, $false, $false, $false, $false, $true, 1, $false, "2+27=", 2) | Out-Null
$d.Content.Find.Execute("66+8=", $true, $false, $false, $false, $false, $true, 1, $false, "52-50=", 2) | Out-Null
$d.Content.Find.Execute("24+35=", $true, $false, $false, $false, $false, $true, 1, $false, "15+33=", 2) | Out-Null
$d.Content.Find.Execute("23+5=", $true, $false, $false, $false, $false, $true, 1, $false, "86-60=", 2) | Out-Null
$d.Content.Find.Execute("97-62=", $true, $false, $false, $false, $false, $true, 1, $false, "51-43=", 2) | Out-Null
$d.Content.Find.Execute("83-25=", $true, $false, $false, $false, $false, $true, 1, $false, "98-74=", 2) | Out-Null
$d.Content.Find.Execute("84-63=", $true, $false, $false, $false, $false, $true, 1, $false, "36+58=", 2) | Out-Null
$d.Content.Find.Execute("38+10=", $true, $false, $false, $false, $false, $true, 1, $false, "0+34=", 2) | Out-Null
$d.Content.Find.Execute("87-45=", $true, $false, $false, $false, $false, $true, 1, $false, "60+32=", 2) | Out-Null
$d.Content.Find.Execute("16+51=", $true, $false, $false, $false, $false, $true, 1, $false, "61-0=", 2) | Out-Null
$d.Content.Find.Execute("83-70=", $true, $false, $false, $false, $false, $true, 1, $false, "55+12=", 2) | Out-Null
$d.Content.Find.Execute("54+10=", $true, $false, $false, $false, $false, $true, 1, $false, "23+53=", 2) | Out-Null
$d.Content.Find.Execute("38+31=", $true, $false, $false, $false, $false, $true, 1, $false, "62-26=", 2) | Out-Null
$d.Content.Find.Execute("8+68=", $true, $false, $false, $false, $false, $true, 1, $false, "37-19=", 2) | Out-Null
$d.Content.Find.Execute("55-38=", $true, $false, $false, $false, $false, $true, 1, $false, "94-57=", 2) | Out-Null
$d.Content.Find.Execute("79-60=", $true, $false, $false, $false, $false, $true, 1, $false, "28+7=", 2) | Out-Null
$d.Content.Find.Execute("79-18=", $true, $false, $false, $false, $false, $true, 1, $false, "2+75=", 2) | Out-Null
$d.Content.Find.Execute("49+32=", $true, $false, $false, $false, $false, $true, 1, $false, "84-14=", 2) | Out-Null
$d.Content.Find.Execute("30+17=", $true, $false, $false, $false, $false, $true, 1, $false, "6+11=", 2) | Out-Null
$d.Content.Find.Execute("84-60=", $true, $false, $false, $false, $false, $true, 1, $false, "30+58=", 2) | Out-Null
$d.Content.Find.Execute("41+24=", $true, $false, $false, $false, $false, $true, 1, $false, "34+1=", 2) | Out-Null
$d.Content.Find.Execute("17+78=", $true, $false, $false, $false, $false, $true, 1, $false, "56-22=", 2) | Out-Null
$d.Content.Find.Execute("34+11=", $true, $false, $false, $false, $false, $true, 1, $false, "91-5=", 2) | Out-Null
$d.Content.Find.Execute("27+71=", $true, $false, $false, $false, $false, $true, 1, $false, "30+23=", 2) | Out-Null
$d.Content.Find.Execute("51+28=", $true, $false, $false, $false, $false, $true, 1, $false, "17+81=", 2) | Out-Null
$d.Content.Find.Execute("95-74=", $true, $false, $false, $false, $false, $true, 1, $false, "66-43=", 2) | Out-Null
$d.Content.Find.Execute("15+71=", $true, $false, $false, $false, $false, $true, 1, $false, "18+69=", 2) | Out-Null
$d.Content.Find.Execute("79-52=", $true, $false, $false, $false, $false, $true, 1, $false, "26+7=", 2) | Out-Null
$d.Content.Find.Execute("10-0=", $true, $false, $false, $false, $false, $true, 1, $false, "2+94=", 2) | Out-Null
$d.Content.Find.Execute("30+34=", $true, $false, $false, $false, $false, $true, 1, $false, "89-18=", 2) | Out-Null
$d.Content.Find.Execute("59-12=", $true, $false, $false, $false, $false, $true, 1, $false, "9+42=", 2) | Out-Null
$d.Content.Find.Execute("35+13=", $true, $false, $false, $false, $false, $true, 1, $false, "36+33=", 2) | Out-Null
$d.Content.Find.Execute("9+84=", $true, $false, $false, $false, $false, $true, 1, $false, "97-23=", 2) | Out-Null
$d.Content.Find.Execute("56+38=", $true, $false, $false, $false, $false, $true, 1, $false, "3+64=", 2) | Out-Null
$d.Content.Find.Execute("52-8=", $true, $false, $false, $false, $false, $true, 1, $false, "28+58=", 2) | Out-Null
$d.Content.Find.Execute("74+0=", $true, $false, $false, $false, $false, $true, 1, $false, "62+34=", 2) | Out-Null
$d.Content.Find.Execute("29+14=", $true, $false, $false, $false, $false, $true, 1, $false, "31+39=", 2) | Out-Null
$d.Content.Find.Execute("41+30=", $true, $false, $false, $false, $false, $true, 1, $false, "70-30=", 2) | Out-Null
$d.Content.Find.Execute("97-71=", $true, $false, $false, $false, $false, $true, 1, $false, "52+4=", 2) | Out-Null
$d.Content.Find.Execute("5+41=", $true, $false, $false, $false, $false, $true, 1, $false, "27+9=", 2) | Out-Null
$d.Content.Find.Execute("35+19=", $true, $false, $false, $false, $false, $true, 1, $false, "62-4=", 2) | Out-Null
$d.Content.Find.Execute("49-14=", $true, $false, $false, $false, $false, $true, 1, $false, "32+41=", 2) | Out-Null
$d.Content.Find.Execute("75-63=", $true, $false, $false, $false, $false, $true, 1, $false, "92-91=", 2) | Out-Null
$d.Content.Find.Execute("39-16=", $true, $false, $false, $false, $false, $true, 1, $false, "16-12=", 2) | Out-Null
$d.Content.Find.Execute("74-36=", $true, $false, $false, $false, $false, $true, 1, $false, "16+44=", 2) | Out-Null
$d.Content.Find.Execute("24+23=", $true, $false, $false, $false, $false, $true, 1, $false, "89-47=", 2) | Out-Null
$d.Content.Find.Execute("39+45=", $true, $false, $false, $false, $false, $true, 1, $false, "76-19=", 2) | Out-Null
$d.Content.Find.Execute("70-43=", $true, $false, $false, $false, $false, $true, 1, $false, "45+52=", 2) | Out-Null
$d.Content.Find.Execute("67-27=", $true, $false, $false, $false, $false, $true, 1, $false, "41-17=", 2) | Out-Null
$d.Content.Find.Execute("69-52=", $true, $false, $false, $false, $false, $true, 1, $false, "94-90=", 2) | Out-Null
$d.Content.Find.Execute("64-37=", $true, $false, $false, $false, $false, $true, 1, $false, "42+51=", 2) | Out-Null
$d.Content.Find.Execute("17+28=", $true, $false, $false, $false, $false, $true, 1, $false, "6+75=", 2) | Out-Null
$d.Content.Find.Execute("29-0=", $true, $false, $false, $false, $false, $true, 1, $false, "60+6=", 2) | Out-Null
$d.Content.Find.Execute("33+63=", $true, $false, $false, $false, $false, $true, 1, $false, "84-48=", 2) | Out-Null
$d.Content.Find.Execute("32-31=", $true, $false, $false, $false, $false, $true, 1, $false, "34+62=", 2) | Out-Null
$d.Content.Find.Execute("56-28=", $true, $false, $false, $false, $false, $true, 1, $false, "71-4=", 2) | Out-Null
$d.Content.Find.Execute("95-54=", $true, $false, $false, $false, $false, $true, 1, $false, "42+42=", 2) | Out-Null
$d.Content.Find.Execute("15+47=", $true, $false, $false, $false, $false, $true, 1, $false, "53+13=", 2) | Out-Null
$d.Content.Find.Execute("25+27=", $true, $false, $false, $false, $false, $true, 1, $false, "40+6=", 2) | Out-Null
$d.Content.Find.Execute("31+48=", $true, $false, $false, $false, $false, $true, 1, $false, "46-32=", 2) | Out-Null
$d.Content.Find.Execute("14+31=", $true, $false, $false, $false, $false, $true, 1, $false, "55-28=", 2) | Out-Null
$d.Content.Find.Execute("83-22=", $true, $false, $false, $false, $false, $true, 1, $false, "28+17=", 2) | Out-Null
$d.Content.Find.Execute("65-64=", $true, $false, $false, $false, $false, $true, 1, $false, "41+15=", 2) | Out-Null
$d.Content.Find.Execute("15+35=", $true, $false, $false, $false, $false, $true, 1, $false, "71-4=", 2) | Out-Null
$d.Content.Find.Execute("21+14=", $true, $false, $false, $false, $false, $true, 1, $false, "49+27=", 2) | Out-Null
$d.Content.Find.Execute("35-15=", $true, $false, $false, $false, $false, $true, 1, $false, "82+12=", 2) | Out-Null
$d.Content.Find.Execute("4+15=", $true, $false, $false, $false, $false, $true, 1, $false, "77+16=", 2) | Out-Null
$d.Content.Find.Execute("54+5=", $true, $false, $false, $false, $false, $true, 1, $false, "54-28=", 2) | Out-Null
$d.Content.Find.Execute("46+38=", $true, $false, $false, $false, $false, $true, 1, $false, "70-49=", 2) | Out-Null
$d.Content.Find.Execute("56-23=", $true, $false, $false, $false, $false, $true, 1, $false, "84-54=", 2) | Out-Null
$d.Content.Find.Execute("38+1=", $true, $false, $false, $false, $false, $true, 1, $false, "29+28=", 2) | Out-Null
$d.Content.Find.Execute("4+9=", $true, $false, $false, $false, $false, $true, 1, $false, "1+65=", 2) | Out-Null
$d.Content.Find.Execute("44+49=", $true, $false, $false, $false, $false, $true, 1, $false, "56+17=", 2) | Out-Null
$d.Content.Find.Execute("72-20=", $true, $false, $false, $false, $false, $true, 1, $false, "59-19=", 2) | Out-Null
$d.Content.Find.Execute("68-14=", $true, $false, $false, $false, $false, $true, 1, $false, "28+20=", 2) | Out-Null
$d.Content.Find.Execute("29+51=", $true, $false, $false, $false, $false, $true, 1, $false, "99-52=", 2) | Out-Null
$d.Content.Find.Execute("75-54=", $true, $false, $false, $false, $false, $true, 1, $false, "77+14=", 2) | Out-Null
$d.Content.Find.Execute("58+36=", $true, $false, $false, $false, $false, $true, 1, $false, "21+78=", 2) | Out-Null
$d.Content.Find.Execute("35-2=", $true, $false, $false, $false, $false, $true, 1, $false, "49-9=", 2) | Out-Null
$d.Content.Find.Execute("48-32=", $true, $false, $false, $false, $false, $true, 1, $false, "11+16=", 2) | Out-Null
